$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-8 (regions, variable label, values, ranking)
$ws.Range("A2").Value = "Rondônia"
$ws.Range("B2").Value = "Variação 2022/2021"
$ws.Range("C2").Value = 25.87011789417062
$ws.Range("D2").Value = "1º"

$ws.Range("A3").Value = "Distrito Federal"
$ws.Range("B3").Value = "Variação 2022/2021"
$ws.Range("C3").Value = 15.2625061432311
$ws.Range("D3").Value = "2º"

$ws.Range("A4").Value = "Mato Grosso"
$ws.Range("B4").Value = "Variação 2022/2021"
$ws.Range("C4").Value = 13.46567778129448
$ws.Range("D4").Value = "3º"

$ws.Range("A5").Value = "Tocantins"
$ws.Range("B5").Value = "Variação 2022/2021"
$ws.Range("C5").Value = 13.3133806028382
$ws.Range("D5").Value = "4º"

$ws.Range("A6").Value = "Goiás"
$ws.Range("B6").Value = "Variação 2022/2021"
$ws.Range("C6").Value = 12.18552997886027
$ws.Range("D6").Value = "5º"

$ws.Range("A7").Value = "Maranhão"
$ws.Range("B7").Value = "Variação 2022/2021"
$ws.Range("C7").Value = 11.60052210713454
$ws.Range("D7").Value = "6º"

$ws.Range("A8").Value = "Sergipe"
$ws.Range("B8").Value = "Variação 2022/2021"
$ws.Range("C8").Value = 6.392991047517027
$ws.Range("D8").Value = "23º"

# Row 9: keep "Nordeste" label, update variable label and value; D9 has no value
$ws.Range("A9").Value = "Nordeste"
$ws.Range("B9").Value = "Variação 2022/2021"
$ws.Range("C9").Value = 8.261124244957617

# Row 10 ("Brasil") is removed entirely
$ws.Rows.Item(10).Delete()

$wb.Save()
